$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.994.88'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.744.05'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.05'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5048'
$ws.Range('E7').Value = '  -5.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2749'
$ws.Range('E8').Value = '  -1.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06188'
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.749.10'
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07249'
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.6512'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.14'
$ws.Range('E13').Value = '  -1.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.681'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.51'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.008.91'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.90'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006911'
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.966.75'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.471'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.741'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('E24').Value = '  +2.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.33'
$ws.Range('E25').Value = '  -2.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.506'
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.31'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.777'
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.02'
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.863'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08199'
$ws.Range('E31').Value = '  -3.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.641'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04662'
$ws.Range('E33').Value = '  +1.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.656'
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9935'
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6095'
$ws.Range('E36').Value = '  -2.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.780'
$ws.Range('E37').Value = '  +2.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01621'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.928'
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.000'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '100.56'
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.3915'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7668'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.996'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1159'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.306'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.67'
$ws.Range('E47').Value = '  +1.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05326'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.70'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.641'
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3457'
$ws.Range('E51').Value = '  -1.16%  '
